$d = $word.ActiveDocument
$t = $d.Tables(1)

$newValues = @(
    "33-27=6",
    "18+53=71",
    "91-26=65",
    "16+76=92",
    "43+38=81",
    "9+62=71",
    "26-19=7",
    "5+89=94",
    "18+3=21",
    "44-26=18",
    "93-59=34",
    "53+9=62",
    "29+58=87",
    "92-75=17",
    "43+49=92",
    "17+68=85",
    "70-13=57",
    "70-11=59",
    "39+56=95",
    "17+37=54",
    "14+9=23",
    "56+8=64",
    "32-24=8",
    "83-44=39",
    "40-36=4",
    "19+5=24",
    "17+24=41",
    "91-77=14",
    "8+6=14",
    "83-37=46",
    "27+44=71",
    "29+36=65",
    "81-44=37",
    "9+19=28",
    "52-5=47",
    "51-2=49",
    "75-47=28",
    "80-31=49",
    "8+7=15",
    "66+16=82",
    "24+59=83",
    "43-35=8",
    "25+7=32",
    "41-23=18",
    "79+16=95",
    "83-27=56",
    "50-22=28",
    "56-28=28",
    "81-59=22",
    "92-13=79",
    "77+5=82",
    "29+29=58",
    "93-18=75",
    "29+15=44",
    "44-18=26",
    "41-12=29",
    "79+5=84",
    "66-57=9",
    "6+59=65",
    "31-7=24",
    "89+6=95",
    "66+8=74",
    "74-59=15",
    "29+44=73",
    "3+9=12",
    "20-18=2",
    "95-57=38",
    "37+35=72",
    "74-67=7",
    "81-78=3",
    "39+36=75",
    "53-19=34",
    "9+57=66",
    "18+44=62",
    "62-13=49",
    "19+57=76",
    "26+17=43",
    "37+18=55",
    "45+6=51",
    "90-16=74",
    "48+14=62",
    "38+45=83",
    "20-13=7",
    "80-38=42",
    "28+37=65",
    "44-29=15",
    "7+36=43",
    "55+29=84",
    "4+19=23",
    "29+26=55",
    "19+74=93",
    "15-8=7",
    "19+78=97",
    "16+27=43",
    "43-35=8",
    "77+7=84",
    "49+44=93",
    "52-23=29",
    "92-4=88",
    "26+8=34"
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx++
    }
}

Write-Host "Updated" $idx "cells"
